# AttributeCrosswalk_Okanogan_EDT.xlsx edit
# Adds a new "Temperature- Adult Holding" row to the HabitatAttribute4 table
# (row 45), extends the table/autofilter to match, adjusts the data
# validation rules accordingly, and updates the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add the new data row (A45:C45) ------------------------------------
# A45 re-uses the existing "Temperature: Daily Maximum" Level 2 Attribute,
# B45 is a brand new RTT Habitat Attribute string, C45 flags it as an
# HQ Pathway Attribute ("x"). Setting these values also creates the new
# shared string entry automatically.
$ws.Range("A45").Value = "Temperature: Daily Maximum"
$ws.Range("B45").Value = "Temperature- Adult Holding"
$ws.Range("C45").Value = "x"

# --- 2. Grow the table / autofilter to include the new row ----------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:B45"))

# --- 3. Fix up data validation -------------------------------------------
# Previously a single (pre-existing / stale) rule covered "C45 B3:B36 B39:B45".
# Drop C45 and B45 from that rule (shrinking it back to B3:B36 / B39:B44),
# then add a dedicated new validation rule on B45 referencing the real list
# range $E$2:$E$28.
$ws.Range("C45").Validation.Delete()
$ws.Range("B45").Validation.Delete()
$ws.Range("B45").Validation.Add(3, 1, 1, "`$E`$2:`$E`$28")

# --- 4. Update the active selection on the sheet ---------------------------
$ws.Range("D44").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 21
$win.ScrollColumn = 1
